$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 30.708183
$ws.Range("H2").Value = 92.124549
$ws.Range("I2").Value = 0.3702499640981372
$ws.Range("J2").Value = 0.3702499640981371
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.499519333333335
$ws.Range("N2").Value = 22.498558
$ws.Range("O2").Value = 0.2063210495448665
$ws.Range("P2").Value = 0.2063210495448665
$ws.Range("Q2").Value = 230.2966121000381
$ws.Range("R2").Value = 2072.669508900342
$ws.Range("S2").Value = 0.0763903611866768
$ws.Range("T2").Value = 0.07639036118667679

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 30.708183
$ws.Range("H3").Value = 92.124549
$ws.Range("I3").Value = 0.3702499640981372
$ws.Range("J3").Value = 0.3702499640981371
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.856814
$ws.Range("N3").Value = 44.570442
$ws.Range("O3").Value = 0.4087293226578609
$ws.Range("P3").Value = 0.4087293226578609
$ws.Range("Q3").Value = 456.225763108962
$ws.Range("R3").Value = 4106.031867980658
$ws.Range("S3").Value = 0.1513320170399289
$ws.Range("T3").Value = 0.1513320170399289

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 30.708183
$ws.Range("H4").Value = 92.124549
$ws.Range("I4").Value = 0.3702499640981372
$ws.Range("J4").Value = 0.3702499640981371
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.992451
$ws.Range("N4").Value = 41.977353
$ws.Range("O4").Value = 0.3849496277972726
$ws.Range("P4").Value = 0.3849496277972725
$ws.Range("Q4").Value = 429.6827459265331
$ws.Range("R4").Value = 3867.144713338797
$ws.Range("S4").Value = 0.1425275858715314
$ws.Range("T4").Value = 0.1425275858715314

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.69729066666667
$ws.Range("H5").Value = 38.091872
$ws.Range("I5").Value = 0.1530918131325759
$ws.Range("J5").Value = 0.1530918131325759
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.499519333333335
$ws.Range("N5").Value = 22.498558
$ws.Range("O5").Value = 0.2063210495448665
$ws.Range("P5").Value = 0.2063210495448665
$ws.Range("Q5").Value = 95.22357683561958
$ws.Range("R5").Value = 857.0121915205762
$ws.Range("S5").Value = 0.03158606356223965
$ws.Range("T5").Value = 0.03158606356223964

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.69729066666667
$ws.Range("H6").Value = 38.091872
$ws.Range("I6").Value = 0.1530918131325759
$ws.Range("J6").Value = 0.1530918131325759
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.856814
$ws.Range("N6").Value = 44.570442
$ws.Range("O6").Value = 0.4087293226578609
$ws.Range("P6").Value = 0.4087293226578609
$ws.Range("Q6").Value = 188.6412857386027
$ws.Range("R6").Value = 1697.771571647424
$ws.Range("S6").Value = 0.06257311308614158
$ws.Range("T6").Value = 0.06257311308614157

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.69729066666667
$ws.Range("H7").Value = 38.091872
$ws.Range("I7").Value = 0.1530918131325759
$ws.Range("J7").Value = 0.1530918131325759
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.992451
$ws.Range("N7").Value = 41.977353
$ws.Range("O7").Value = 0.3849496277972726
$ws.Range("P7").Value = 0.3849496277972725
$ws.Range("Q7").Value = 177.6662174860907
$ws.Range("R7").Value = 1598.995957374816
$ws.Range("S7").Value = 0.05893263648419472
$ws.Range("T7").Value = 0.0589326364841947

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 39.53358366666667
$ws.Range("H8").Value = 118.600751
$ws.Range("I8").Value = 0.476658222769287
$ws.Range("J8").Value = 0.476658222769287
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.499519333333335
$ws.Range("N8").Value = 22.498558
$ws.Range("O8").Value = 0.2063210495448665
$ws.Range("P8").Value = 0.2063210495448665
$ws.Range("Q8").Value = 296.4828750241176
$ws.Range("R8").Value = 2668.345875217058
$ws.Range("S8").Value = 0.09834462479595009
$ws.Range("T8").Value = 0.09834462479595008

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 39.53358366666667
$ws.Range("H9").Value = 118.600751
$ws.Range("I9").Value = 0.476658222769287
$ws.Range("J9").Value = 0.476658222769287
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 14.856814
$ws.Range("N9").Value = 44.570442
$ws.Range("O9").Value = 0.4087293226578609
$ws.Range("P9").Value = 0.4087293226578609
$ws.Range("Q9").Value = 587.3430992891047
$ws.Range("R9").Value = 5286.087893601942
$ws.Range("S9").Value = 0.1948241925317905
$ws.Range("T9").Value = 0.1948241925317904

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 39.53358366666667
$ws.Range("H10").Value = 118.600751
$ws.Range("I10").Value = 0.476658222769287
$ws.Range("J10").Value = 0.476658222769287
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.992451
$ws.Range("N10").Value = 41.977353
$ws.Range("O10").Value = 0.3849496277972726
$ws.Range("P10").Value = 0.3849496277972725
$ws.Range("Q10").Value = 553.1717323102337
$ws.Range("R10").Value = 4978.545590792103
$ws.Range("S10").Value = 0.1834894054415465
$ws.Range("T10").Value = 0.1834894054415464
